$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number but must remain stored as
# TEXT (matching the source data which uses "." as a thousands separator for
# other rows). Force text format before assigning so Excel does not
# auto-convert the string into a numeric value.
$textCells = @("D4", "D5", "D6", "D10", "D14", "D15", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D29", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D44", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "60.823.10"
$ws.Range("E2").Value = "  -1.41%  "

$ws.Range("D3").Value = "3.401.25"
$ws.Range("E3").Value = "  -1.36%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "572.72"
$ws.Range("E5").Value = "  -0.93%  "

$ws.Range("D6").Value = "142.70"
$ws.Range("E6").Value = "  -3.09%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").Value = "3.404.31"
$ws.Range("E8").Value = "  -1.28%  "

$ws.Range("E9").Value = "  +0.45%  "

$ws.Range("D10").Value = "7.56"
$ws.Range("E10").Value = "  -1.46%  "

$ws.Range("E11").Value = "  +1.21%  "

$ws.Range("E12").Value = "  +2.30%  "

$ws.Range("D13").Value = "3.984.61"
$ws.Range("E13").Value = "  -1.18%  "

$ws.Range("D14").Value = "28.32"

$ws.Range("D15").Value = "0.124"
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("E16").Value = "  -1.16%  "

$ws.Range("D17").Value = "3.409.35"
$ws.Range("E17").Value = "  -1.02%  "

$ws.Range("D18").Value = "60.963.77"
$ws.Range("E18").Value = "  -1.19%  "

$ws.Range("D19").Value = "6.33"
$ws.Range("E19").Value = "  +0.88%  "

$ws.Range("D20").Value = "14.24"
$ws.Range("E20").Value = "  +1.29%  "

$ws.Range("D21").Value = "9.21"
$ws.Range("E21").Value = "  -1.70%  "

$ws.Range("D22").Value = "389.19"
$ws.Range("E22").Value = "  +1.65%  "

$ws.Range("D23").Value = "0.566"
$ws.Range("E23").Value = "  +0.24%  "

$ws.Range("D24").Value = "73.00"
$ws.Range("E24").Value = "  +1.08%  "

$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.13%  "

$ws.Range("D26").Value = "0.0000122"
$ws.Range("E26").Value = "  -0.93%  "

$ws.Range("D27").Value = "3.547.69"
$ws.Range("E27").Value = "  -1.08%  "

$ws.Range("E28").Value = "  +0.29%  "

$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("E30").Value = "  -3.98%  "

$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("E32").Value = "  -6.49%  "

$ws.Range("E33").Value = "  +0.15%  "

$ws.Range("E34").Value = "  -0.07%  "

$ws.Range("D35").Value = "23.83"
$ws.Range("E35").Value = "  -0.68%  "

$ws.Range("D36").Value = "7.03"
$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("D37").Value = "3.431.63"
$ws.Range("E37").Value = "  -1.02%  "

$ws.Range("D38").Value = "5.11"
$ws.Range("E38").Value = "  -1.67%  "

$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "1.54"
$ws.Range("E39").Value = "  -0.55%  "

$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "167.14"
$ws.Range("E40").Value = "  +0.85%  "

$ws.Range("D41").Value = "0.0786"
$ws.Range("E41").Value = "  +0.52%  "

$ws.Range("D42").Value = "27.17"
$ws.Range("E42").Value = "  +6.32%  "

$ws.Range("E43").Value = "  -0.13%  "

$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("E45").Value = "  +0.57%  "

$ws.Range("E46").Value = "  -0.53%  "

$ws.Range("E47").Value = "  -1.48%  "

$ws.Range("D48").Value = "2.573.76"
$ws.Range("E48").Value = "  -1.25%  "

$ws.Range("E49").Value = "  -3.45%  "

$ws.Range("E50").Value = "  +1.05%  "

$ws.Range("D51").Value = "23.13"
$ws.Range("E51").Value = "  -1.63%  "

# Remove the temporary text-number-format override so the cell style
# matches the original workbook (no explicit style index) while keeping
# the stored value as text.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
